$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6607595682144165
$ws.Range("B1").Value = 1.01718008518219
$ws.Range("C1").Value = 2.324509382247925
$ws.Range("D1").Value = 3.929444313049316
$ws.Range("E1").Value = 1.567742109298706
